# Updated customer delete logic and fixed typos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Delete old customer record (CAN 1001, Rahul Sharma) ---
# Remaining rows 3..6 shift up to become rows 2..5.
$ws.Rows.Item(2).Delete()

# --- Row 2 (was row 3): Priya Verma ---
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "1002"
$ws.Cells.Item(2, 4).NumberFormat = "General"
$ws.Cells.Item(2, 4).Value = 9123456789
$ws.Cells.Item(2, 7).NumberFormat = "General"
$ws.Cells.Item(2, 7).Value = 650

# --- Row 3 (was row 4): Amitabh Patel ---
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "1003"
$ws.Cells.Item(3, 4).NumberFormat = "General"
$ws.Cells.Item(3, 4).Value = 9988776655
$ws.Cells.Item(3, 7).NumberFormat = "General"
$ws.Cells.Item(3, 7).Value = 400

# --- Row 4 (was row 5): Sneha Gupta ---
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "1004"
$ws.Cells.Item(4, 4).NumberFormat = "General"
$ws.Cells.Item(4, 4).Value = 9000011111
$ws.Cells.Item(4, 7).NumberFormat = "General"
$ws.Cells.Item(4, 7).Value = 1000

# --- Row 5 (was row 6): Vikram Singh ---
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "1005"
$ws.Cells.Item(5, 4).NumberFormat = "General"
$ws.Cells.Item(5, 4).Value = 8888822222
$ws.Cells.Item(5, 7).NumberFormat = "General"
$ws.Cells.Item(5, 7).Value = 500

# --- Add two new customer records, row 6 (typo fixes) and rows 7-8 (new) ---

# Row 6: rahul verma / nagpore
$ws.Cells.Item(6, 1).NumberFormat = "@"
$ws.Cells.Item(6, 1).Value = "6654"
$ws.Cells.Item(6, 2).Value = "rahul verma"
$ws.Cells.Item(6, 3).Value = "nagpore"
$ws.Cells.Item(6, 4).Value = 645451
$ws.Cells.Item(6, 5).Value = 616546546
$ws.Cells.Item(6, 6).NumberFormat = "@"
$ws.Cells.Item(6, 6).Value = "2025-12-20"
$ws.Cells.Item(6, 7).Value = 200

# Row 7: Rahul Kumar
$ws.Cells.Item(7, 1).NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = "65541654"
$ws.Cells.Item(7, 2).Value = "Rahul Kumar"
$ws.Cells.Item(7, 3).Value = "12 koli"
$ws.Cells.Item(7, 4).Value = 264768446
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = "68146545314661"
$ws.Cells.Item(7, 6).Value = ""
$ws.Cells.Item(7, 7).Value = 0

# Row 8: Rahul Sharma
$ws.Cells.Item(8, 1).NumberFormat = "@"
$ws.Cells.Item(8, 1).Value = "23436854"
$ws.Cells.Item(8, 2).Value = "Rahul Sharma"
$ws.Cells.Item(8, 3).Value = "10 no puliya"
$ws.Cells.Item(8, 4).Value = 987845665
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = "234143434788"
$ws.Cells.Item(8, 6).NumberFormat = "@"
$ws.Cells.Item(8, 6).Value = "2025-12-15"
$ws.Cells.Item(8, 7).NumberFormat = "@"
$ws.Cells.Item(8, 7).Value = "1"
